$wb = $excel.ActiveWorkbook

# ---- Sheet 1: Means ----
$ws1 = $wb.Worksheets.Item("Means")

# New header columns F and G
$ws1.Range("F1").Value = "Within 5 miles of HFC production facility"
$ws1.Range("G1").Value = "Within 10 miles of HFC production facility"

# New data columns F and G for rows 2-10
$meansF = @(56, 39, 5.3, 4.7, 79, 6.3, 8.3, 93, 0.51)
$meansG = @(68, 27, 5.7, 5, 80, 5.3, 5.4, 81, 0.5)
for ($i = 0; $i -lt 9; $i++) {
    $row = $i + 2
    $ws1.Cells.Item($row, 6).Value = $meansF[$i]
    $ws1.Cells.Item($row, 7).Value = $meansG[$i]
}

# Updated values in row 9 (Total Cancer Risk) and row 10 (Total Respiratory)
$ws1.Range("B9").Value = 26
$ws1.Range("C9").Value = 39
$ws1.Range("D9").Value = 93
$ws1.Range("E9").Value = 95

$ws1.Range("B10").Value = 0.32
$ws1.Range("C10").Value = 0.43
$ws1.Range("D10").Value = 0.5
$ws1.Range("E10").Value = 0.5

# ---- Sheet 2: Standard Deviations ----
$ws2 = $wb.Worksheets.Item("Standard Deviations")

$ws2.Range("F1").Value = "Within 5 mile of HFC production facility SD"
$ws2.Range("G1").Value = "Within 10 mile of HFC production facility SD"

$sdF = @(23, 23, 8.1, 9, 34, 7.9, 11, 32, 0.039)
$sdG = @(29, 29, 5.8, 6.2, 31, 7.1, 8.6, 35, 0.065)
for ($i = 0; $i -lt 9; $i++) {
    $row = $i + 2
    $ws2.Cells.Item($row, 6).Value = $sdF[$i]
    $ws2.Cells.Item($row, 7).Value = $sdG[$i]
}

$ws2.Range("B9").Value = 8.6
$ws2.Range("C9").Value = 24
$ws2.Range("D9").Value = 14
$ws2.Range("E9").Value = 15

$ws2.Range("B10").Value = 0.14
$ws2.Range("C10").Value = 0.084
$ws2.Range("D10").Value = 0
$ws2.Range("E10").Value = 0
